$wb = $excel.ActiveWorkbook

# "Generate Report for Archive": the handoff status moves from
# "Ready for handoff" to "In Translation" everywhere it appears, and the
# status column on each sheet is narrowed to match the new, shorter label.

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

foreach ($sheet in @($overview, $zhcn, $dede)) {
    $sheet.Cells.Replace("Ready for handoff", "In Translation")
}

# Narrow the now-shorter status columns (characters-width units on this
# host are quantized to 1/6ths; 12.5 lands in the middle of the bucket
# that serializes to the target raw width).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
